$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived values (NATMI LR-pairs output) for Vcan-Cd44 sheet.
# Columns: G=Ligand avg expr, H=Ligand total expr, I/J=Ligand derived specificity (avg/total),
#          M=Receptor avg expr, N=Receptor total expr, O/P=Receptor derived specificity (avg/total),
#          Q=Edge avg expr weight, R=Edge total expr weight, S/T=Edge derived specificity (avg/total)

$ws.Cells.Item(2, 7).Value = 3.262296333333333
$ws.Cells.Item(2, 8).Value = 9.786889
$ws.Cells.Item(2, 9).Value = 0.01915820289899999
$ws.Cells.Item(2, 10).Value = 0.01915820289899999
$ws.Cells.Item(2, 13).Value = 16.27546433333333
$ws.Cells.Item(2, 14).Value = 48.826393
$ws.Cells.Item(2, 15).Value = 0.06628560529319844
$ws.Cells.Item(2, 16).Value = 0.06628560529319844
$ws.Cells.Item(2, 17).Value = 53.09538761793078
$ws.Cells.Item(2, 18).Value = 477.858488561377
$ws.Cells.Item(2, 19).Value = 0.001269913075490123
$ws.Cells.Item(2, 20).Value = 0.001269913075490123
$ws.Cells.Item(3, 7).Value = 3.262296333333333
$ws.Cells.Item(3, 8).Value = 9.786889
$ws.Cells.Item(3, 9).Value = 0.01915820289899999
$ws.Cells.Item(3, 10).Value = 0.01915820289899999
$ws.Cells.Item(3, 15).Value = 0.3480686258826592
$ws.Cells.Item(3, 16).Value = 0.3480686258826592
$ws.Cells.Item(3, 17).Value = 278.8062133118461
$ws.Cells.Item(3, 18).Value = 2509.255919806615
$ws.Cells.Item(3, 19).Value = 0.006668369357436105
$ws.Cells.Item(3, 20).Value = 0.006668369357436105
$ws.Cells.Item(4, 7).Value = 3.262296333333333
$ws.Cells.Item(4, 8).Value = 9.786889
$ws.Cells.Item(4, 9).Value = 0.01915820289899999
$ws.Cells.Item(4, 10).Value = 0.01915820289899999
$ws.Cells.Item(4, 13).Value = 42.61351133333333
$ws.Cells.Item(4, 14).Value = 127.840534
$ws.Cells.Item(4, 15).Value = 0.17355341356458
$ws.Cells.Item(4, 16).Value = 0.17355341356458
$ws.Cells.Item(4, 17).Value = 139.0179017731918
$ws.Cells.Item(4, 18).Value = 1251.161115958726
$ws.Cells.Item(4, 19).Value = 0.003324971510884281
$ws.Cells.Item(4, 20).Value = 0.003324971510884281
$ws.Cells.Item(5, 7).Value = 3.262296333333333
$ws.Cells.Item(5, 8).Value = 9.786889
$ws.Cells.Item(5, 9).Value = 0.01915820289899999
$ws.Cells.Item(5, 10).Value = 0.01915820289899999
$ws.Cells.Item(5, 13).Value = 101.183272
$ws.Cells.Item(5, 14).Value = 303.549816
$ws.Cells.Item(5, 15).Value = 0.4120923552595624
$ws.Cells.Item(5, 16).Value = 0.4120923552595624
$ws.Cells.Item(5, 17).Value = 330.0898172402693
$ws.Cells.Item(5, 18).Value = 2970.808355162424
$ws.Cells.Item(5, 19).Value = 0.007894948955189482
$ws.Cells.Item(5, 20).Value = 0.007894948955189482
$ws.Cells.Item(6, 9).Value = 0.8527862647199704
$ws.Cells.Item(6, 10).Value = 0.8527862647199704
$ws.Cells.Item(6, 13).Value = 16.27546433333333
$ws.Cells.Item(6, 14).Value = 48.826393
$ws.Cells.Item(6, 15).Value = 0.06628560529319844
$ws.Cells.Item(6, 16).Value = 0.06628560529319844
$ws.Cells.Item(6, 17).Value = 2363.427171079684
$ws.Cells.Item(6, 18).Value = 21270.84453971716
$ws.Cells.Item(6, 19).Value = 0.05652745374268899
$ws.Cells.Item(6, 20).Value = 0.05652745374268899
$ws.Cells.Item(7, 9).Value = 0.8527862647199704
$ws.Cells.Item(7, 10).Value = 0.8527862647199704
$ws.Cells.Item(7, 15).Value = 0.3480686258826592
$ws.Cells.Item(7, 16).Value = 0.3480686258826592
$ws.Cells.Item(7, 19).Value = 0.2968281433326858
$ws.Cells.Item(7, 20).Value = 0.2968281433326858
$ws.Cells.Item(8, 9).Value = 0.8527862647199704
$ws.Cells.Item(8, 10).Value = 0.8527862647199704
$ws.Cells.Item(8, 13).Value = 42.61351133333333
$ws.Cells.Item(8, 14).Value = 127.840534
$ws.Cells.Item(8, 15).Value = 0.17355341356458
$ws.Cells.Item(8, 16).Value = 0.17355341356458
$ws.Cells.Item(8, 17).Value = 6188.083392130485
$ws.Cells.Item(8, 18).Value = 55692.75052917437
$ws.Cells.Item(8, 19).Value = 0.1480039672831384
$ws.Cells.Item(8, 20).Value = 0.1480039672831384
$ws.Cells.Item(9, 9).Value = 0.8527862647199704
$ws.Cells.Item(9, 10).Value = 0.8527862647199704
$ws.Cells.Item(9, 13).Value = 101.183272
$ws.Cells.Item(9, 14).Value = 303.549816
$ws.Cells.Item(9, 15).Value = 0.4120923552595624
$ws.Cells.Item(9, 16).Value = 0.4120923552595624
$ws.Cells.Item(9, 17).Value = 14693.23943119531
$ws.Cells.Item(9, 18).Value = 132239.1548807578
$ws.Cells.Item(9, 19).Value = 0.3514267003614573
$ws.Cells.Item(9, 20).Value = 0.3514267003614573
$ws.Cells.Item(10, 7).Value = 21.305189
$ws.Cells.Item(10, 8).Value = 63.915567
$ws.Cells.Item(10, 9).Value = 0.1251171236325075
$ws.Cells.Item(10, 10).Value = 0.1251171236325075
$ws.Cells.Item(10, 13).Value = 16.27546433333333
$ws.Cells.Item(10, 14).Value = 48.826393
$ws.Cells.Item(10, 15).Value = 0.06628560529319844
$ws.Cells.Item(10, 16).Value = 0.06628560529319844
$ws.Cells.Item(10, 17).Value = 346.7518436844256
$ws.Cells.Item(10, 18).Value = 3120.76659315983
$ws.Cells.Item(10, 19).Value = 0.008293464272524705
$ws.Cells.Item(10, 20).Value = 0.008293464272524705
$ws.Cells.Item(11, 7).Value = 21.305189
$ws.Cells.Item(11, 8).Value = 63.915567
$ws.Cells.Item(11, 9).Value = 0.1251171236325075
$ws.Cells.Item(11, 10).Value = 0.1251171236325075
$ws.Cells.Item(11, 15).Value = 0.3480686258826592
$ws.Cells.Item(11, 16).Value = 0.3480686258826592
$ws.Cells.Item(11, 17).Value = 1820.809166932372
$ws.Cells.Item(11, 18).Value = 16387.28250239135
$ws.Cells.Item(11, 19).Value = 0.04354934529715768
$ws.Cells.Item(11, 20).Value = 0.04354934529715768
$ws.Cells.Item(12, 7).Value = 21.305189
$ws.Cells.Item(12, 8).Value = 63.915567
$ws.Cells.Item(12, 9).Value = 0.1251171236325075
$ws.Cells.Item(12, 10).Value = 0.1251171236325075
$ws.Cells.Item(12, 13).Value = 42.61351133333333
$ws.Cells.Item(12, 14).Value = 127.840534
$ws.Cells.Item(12, 15).Value = 0.17355341356458
$ws.Cells.Item(12, 16).Value = 0.17355341356458
$ws.Cells.Item(12, 17).Value = 907.8889129103085
$ws.Cells.Item(12, 18).Value = 8171.000216192777
$ws.Cells.Item(12, 19).Value = 0.02171450390180327
$ws.Cells.Item(12, 20).Value = 0.02171450390180327
$ws.Cells.Item(13, 7).Value = 21.305189
$ws.Cells.Item(13, 8).Value = 63.915567
$ws.Cells.Item(13, 9).Value = 0.1251171236325075
$ws.Cells.Item(13, 10).Value = 0.1251171236325075
$ws.Cells.Item(13, 13).Value = 101.183272
$ws.Cells.Item(13, 14).Value = 303.549816
$ws.Cells.Item(13, 15).Value = 0.4120923552595624
$ws.Cells.Item(13, 16).Value = 0.4120923552595624
$ws.Cells.Item(13, 17).Value = 2155.728733598407
$ws.Cells.Item(13, 18).Value = 19401.55860238567
$ws.Cells.Item(13, 19).Value = 0.05155981016102189
$ws.Cells.Item(13, 20).Value = 0.05155981016102189
$ws.Cells.Item(14, 7).Value = 0.500358
$ws.Cells.Item(14, 8).Value = 1.501074
$ws.Cells.Item(14, 9).Value = 0.002938408748521978
$ws.Cells.Item(14, 10).Value = 0.002938408748521978
$ws.Cells.Item(14, 13).Value = 16.27546433333333
$ws.Cells.Item(14, 14).Value = 48.826393
$ws.Cells.Item(14, 15).Value = 0.06628560529319844
$ws.Cells.Item(14, 16).Value = 0.06628560529319844
$ws.Cells.Item(14, 17).Value = 8.143558782897999
$ws.Cells.Item(14, 18).Value = 73.29202904608199
$ws.Cells.Item(14, 19).Value = 0.000194774202494609
$ws.Cells.Item(14, 20).Value = 0.000194774202494609
$ws.Cells.Item(15, 7).Value = 0.500358
$ws.Cells.Item(15, 8).Value = 1.501074
$ws.Cells.Item(15, 9).Value = 0.002938408748521978
$ws.Cells.Item(15, 10).Value = 0.002938408748521978
$ws.Cells.Item(15, 15).Value = 0.3480686258826592
$ws.Cells.Item(15, 16).Value = 0.3480686258826592
$ws.Cells.Item(15, 17).Value = 42.76218498451
$ws.Cells.Item(15, 18).Value = 384.85966486059
$ws.Cells.Item(15, 19).Value = 0.001022767895379629
$ws.Cells.Item(15, 20).Value = 0.001022767895379629
$ws.Cells.Item(16, 7).Value = 0.500358
$ws.Cells.Item(16, 8).Value = 1.501074
$ws.Cells.Item(16, 9).Value = 0.002938408748521978
$ws.Cells.Item(16, 10).Value = 0.002938408748521978
$ws.Cells.Item(16, 13).Value = 42.61351133333333
$ws.Cells.Item(16, 14).Value = 127.840534
$ws.Cells.Item(16, 15).Value = 0.17355341356458
$ws.Cells.Item(16, 16).Value = 0.17355341356458
$ws.Cells.Item(16, 17).Value = 21.322011303724
$ws.Cells.Item(16, 18).Value = 191.898101733516
$ws.Cells.Item(16, 19).Value = 0.0005099708687540147
$ws.Cells.Item(16, 20).Value = 0.0005099708687540148
$ws.Cells.Item(17, 7).Value = 0.500358
$ws.Cells.Item(17, 8).Value = 1.501074
$ws.Cells.Item(17, 9).Value = 0.002938408748521978
$ws.Cells.Item(17, 10).Value = 0.002938408748521978
$ws.Cells.Item(17, 13).Value = 101.183272
$ws.Cells.Item(17, 14).Value = 303.549816
$ws.Cells.Item(17, 15).Value = 0.4120923552595624
$ws.Cells.Item(17, 16).Value = 0.4120923552595624
$ws.Cells.Item(17, 17).Value = 50.62785961137599
$ws.Cells.Item(17, 18).Value = 455.6507365023839
$ws.Cells.Item(17, 19).Value = 0.001210895781893725
$ws.Cells.Item(17, 20).Value = 0.001210895781893725
